$wb = $excel.ActiveWorkbook

# --- Sheet "National level": update row 2 values ---
$ws1 = $wb.Worksheets.Item("National level")
$nat = New-Object 'object[,]' 1,5
$nat[0,0] = 0
$nat[0,1] = 0.1838134430727023
$nat[0,2] = 0.1906721536351166
$nat[0,3] = 0.1920438957475995
$nat[0,4] = 0.1810699588477366
$ws1.Range("A2:E2").Value = $nat

# --- Sheet "State level": update full data table A2:F33 ---
$ws2 = $wb.Worksheets.Item("State level")
$state = New-Object 'object[,]' 32,6
$state[0,0] = "Bihar"
$state[0,1] = 0
$state[0,2] = 0.02631578947368421
$state[0,3] = 0.07894736842105263
$state[0,4] = 0.1578947368421053
$state[0,5] = 0.7368421052631579
$state[1,0] = "Jharkhand"
$state[1,1] = 0
$state[1,2] = 0.08333333333333333
$state[1,3] = 0.08333333333333333
$state[1,4] = 0.3333333333333333
$state[1,5] = 0.5
$state[2,0] = "Puducherry"
$state[2,1] = 0
$state[2,2] = 0
$state[2,3] = 0
$state[2,4] = 0.25
$state[2,5] = 0.5
$state[3,0] = "Daman and Diu"
$state[3,1] = 0
$state[3,2] = 0
$state[3,3] = 0.5
$state[3,4] = 0
$state[3,5] = 0.5
$state[4,0] = "Mizoram"
$state[4,1] = 0
$state[4,2] = 0.1818181818181818
$state[4,3] = 0.1818181818181818
$state[4,4] = 0.1818181818181818
$state[4,5] = 0.4545454545454545
$state[5,0] = "Uttar Pradesh"
$state[5,1] = 0
$state[5,2] = 0.06666666666666667
$state[5,3] = 0.1866666666666667
$state[5,4] = 0.28
$state[5,5] = 0.4533333333333333
$state[6,0] = "Arunachal Pradesh"
$state[6,1] = 0
$state[6,2] = 0.12
$state[6,3] = 0.2
$state[6,4] = 0.04
$state[6,5] = 0.4
$state[7,0] = "Nagaland"
$state[7,1] = 0
$state[7,2] = 0.3636363636363637
$state[7,3] = 0.1818181818181818
$state[7,4] = 0
$state[7,5] = 0.3636363636363637
$state[8,0] = "Manipur"
$state[8,1] = 0
$state[8,2] = 0.0625
$state[8,3] = 0.25
$state[8,4] = 0.3125
$state[8,5] = 0.3125
$state[9,0] = "Punjab"
$state[9,1] = 0
$state[9,2] = 0.1363636363636364
$state[9,3] = 0.1818181818181818
$state[9,4] = 0.4090909090909091
$state[9,5] = 0.2727272727272727
$state[10,0] = "Telangana"
$state[10,1] = 0
$state[10,2] = 0.2424242424242424
$state[10,3] = 0.2727272727272727
$state[10,4] = 0.2121212121212121
$state[10,5] = 0.1515151515151515
$state[11,0] = "Odisha"
$state[11,1] = 0
$state[11,2] = 0.2333333333333333
$state[11,3] = 0.2
$state[11,4] = 0.3666666666666666
$state[11,5] = 0.1333333333333333
$state[12,0] = "Tripura"
$state[12,1] = 0
$state[12,2] = 0
$state[12,3] = 0.375
$state[12,4] = 0.5
$state[12,5] = 0.125
$state[13,0] = "Meghalaya"
$state[13,1] = 0
$state[13,2] = 0.2727272727272727
$state[13,3] = 0
$state[13,4] = 0.5454545454545454
$state[13,5] = 0.09090909090909091
$state[14,0] = "Jammu and Kashmir"
$state[14,1] = 0
$state[14,2] = 0.1818181818181818
$state[14,3] = 0.3181818181818182
$state[14,4] = 0.3636363636363637
$state[14,5] = 0.09090909090909091
$state[15,0] = "Delhi"
$state[15,1] = 0
$state[15,2] = 0.2727272727272727
$state[15,3] = 0.2727272727272727
$state[15,4] = 0.3636363636363637
$state[15,5] = 0.09090909090909091
$state[16,0] = "Haryana"
$state[16,1] = 0
$state[16,2] = 0.1363636363636364
$state[16,3] = 0.4545454545454545
$state[16,4] = 0.2727272727272727
$state[16,5] = 0.09090909090909091
$state[17,0] = "Assam"
$state[17,1] = 0
$state[17,2] = 0.1818181818181818
$state[17,3] = 0.1818181818181818
$state[17,4] = 0.1818181818181818
$state[17,5] = 0.09090909090909091
$state[18,0] = "Madhya Pradesh"
$state[18,1] = 0
$state[18,2] = 0.1538461538461539
$state[18,3] = 0.3076923076923077
$state[18,4] = 0.3269230769230769
$state[18,5] = 0.0576923076923077
$state[19,0] = "West Bengal"
$state[19,1] = 0
$state[19,2] = 0.391304347826087
$state[19,3] = 0.2173913043478261
$state[19,4] = 0.08695652173913043
$state[19,5] = 0.04347826086956522
$state[20,0] = "Chhattisgarh"
$state[20,1] = 0
$state[20,2] = 0.2592592592592592
$state[20,3] = 0.4074074074074074
$state[20,4] = 0.1111111111111111
$state[20,5] = 0.03703703703703704
$state[21,0] = "Rajasthan"
$state[21,1] = 0
$state[21,2] = 0.2727272727272727
$state[21,3] = 0.2727272727272727
$state[21,4] = 0.1818181818181818
$state[21,5] = 0.0303030303030303
$state[22,0] = "Himachal Pradesh"
$state[22,1] = 0
$state[22,2] = 0.3333333333333333
$state[22,3] = 0.1666666666666667
$state[22,4] = 0.1666666666666667
$state[22,5] = 0
$state[23,0] = "Gujarat"
$state[23,1] = 0
$state[23,2] = 0.3939393939393939
$state[23,3] = 0.09090909090909091
$state[23,4] = 0.06060606060606061
$state[23,5] = 0
$state[24,0] = "Maharashtra"
$state[24,1] = 0
$state[24,2] = 0.1944444444444445
$state[24,3] = 0.1944444444444445
$state[24,4] = 0.05555555555555555
$state[24,5] = 0
$state[25,0] = "Karnataka"
$state[25,1] = 0
$state[25,2] = 0.2666666666666667
$state[25,3] = 0.03333333333333333
$state[25,4] = 0.03333333333333333
$state[25,5] = 0
$state[26,0] = "Tamil Nadu"
$state[26,1] = 0
$state[26,2] = 0.2702702702702703
$state[26,3] = 0.02702702702702703
$state[26,4] = 0.02702702702702703
$state[26,5] = 0
$state[27,0] = "Chandigarh"
$state[27,1] = 0
$state[27,2] = 0
$state[27,3] = 1
$state[27,4] = 0
$state[27,5] = 0
$state[28,0] = "Dadra and Nagar Haveli"
$state[28,1] = 0
$state[28,2] = 0
$state[28,3] = 1
$state[28,4] = 0
$state[28,5] = 0
$state[29,0] = "Andhra Pradesh"
$state[29,1] = 0
$state[29,2] = 0.3076923076923077
$state[29,3] = 0.1538461538461539
$state[29,4] = 0
$state[29,5] = 0
$state[30,0] = "Uttarakhand"
$state[30,1] = 0
$state[30,2] = 0
$state[30,3] = 0.07692307692307693
$state[30,4] = 0
$state[30,5] = 0
$state[31,0] = "Ladakh"
$state[31,1] = 0
$state[31,2] = 0.5
$state[31,3] = 0
$state[31,4] = 0
$state[31,5] = 0
$ws2.Range("A2:F33").Value = $state

# Remove now-obsolete rows 34 and 35 (Sikkim, Kerala)
$ws2.Rows.Item(35).Delete() | Out-Null
$ws2.Rows.Item(34).Delete() | Out-Null
